$wb = $excel.ActiveWorkbook

# Add the new "empty ER sheet" as the last sheet in the workbook (METABOLIGHTS_METABOLOMICS)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "METABOLIGHTS_METABOLOMICS"

# Header row + data rows for the ER (external reference) sheet
$data = @(
    @("", "TermSourceRef", "Ontology", "TAN", "Content type (validation)", "Notes during templating", "Target term", "Instruction", "Requirement (m/o/n)", "Value (cv/s/d)", "Additional information", "Review comments"),
    @("Source Name", "", "", "", "", "", "", "", "", "", "", ""),
    @("Sample Name", "", "", "", "", "", "", "", "", "", "", ""),
    @("Data File Name", "", "", "", "", "", "", "", "", "", "", ""),
    @("Parameter [area normalization]", "MS:1001999", "MS", "http://purl.obolibrary.org/obo/MS_1001999", "", "", "", "", "", "", "", ""),
    @("Parameter [data transformation]", "MS:1000452", "MS", "http://purl.obolibrary.org/obo/MS_1000452", "", "", "", "", "", "", "", ""),
    @("Parameter [processed data file]", "MS:1003084", "MS", "http://purl.obolibrary.org/obo/MS_1003084", "", "", "", "", "", "", "", ""),
    @("Parameter [Metabolite Assignment File]", "NFDI4PSO:0000077", "NFDI4PSO", "http://purl.obolibrary.org/obo/NFDI4PSO_0000077", "", "", "", "", "", "", "", "")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Column widths (best-fit approximations)
$widths = @(37.28515625, 17.5703125, 9.85546875, 46.5703125, 23.5703125, 22.85546875, 11.28515625, 10.5703125, 20.42578125, 13.5703125, 21.5703125, 17.42578125)
for ($c = 0; $c -lt $widths.Length; $c++) {
    $ws.Columns.Item($c + 1).ColumnWidth = $widths[$c]
}

# Page margins (top/bottom = 2cm, matching the rest of the workbook's sheets)
$ws.PageSetup.TopMargin = 56.692913399999995
$ws.PageSetup.BottomMargin = 56.692913399999995

# Select the full rows (matches the sqref selection state saved on the sheet)
$ws.Range("A1:XFD1048576").Select() | Out-Null
Write-Host "Added METABOLIGHTS_METABOLOMICS sheet with data"
